$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Driver Summary")

# Row 3: Intel(R) Wireless-AC 9560 160MHz - 22.80.0.9 -> Intel(R) Wi-Fi 6 AX201 160MHz - 23.30.0.6
$ws.Range("A3").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 23.30.0.6"
$ws.Range("C3").Value = 83
$ws.Range("D3").Value = 97.7

# Row 4: Intel(R) Wi-Fi 6 AX201 160MHz - 23.30.0.6 -> Intel(R) Wi-Fi 6 AX201 160MHz - 23.50.0.6
$ws.Range("A4").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 23.50.0.6"
$ws.Range("C4").Value = 244
$ws.Range("D4").Value = 98

# Row 5: Intel(R) Dual Band Wireless-AC 3165 - 19.51.40.1 (unchanged name)
$ws.Range("C5").Value = 1254

# Row 6: Intel(R) Wi-Fi 6 AX201 160MHz - 23.50.0.6 -> Intel(R) Wi-Fi 6 AX201 160MHz - 22.150.0.3
$ws.Range("A6").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 22.150.0.3"
$ws.Range("C6").Value = 40
$ws.Range("D6").Value = 98.5

# Row 7: Intel(R) Wi-Fi 6 AX201 160MHz - 22.150.0.3 -> Intel(R) Wireless-AC 9560 160MHz - 22.80.0.9
$ws.Range("A7").Value = "Intel(R) Wireless-AC 9560 160MHz - 22.80.0.9"
$ws.Range("C7").Value = 3
$ws.Range("D7").Value = 98.8

# Row 8: Intel(R) Dual Band Wireless-AC 8265 - 20.90.0.100 (unchanged name)
$ws.Range("C8").Value = 158

# Row 9: Totals
$ws.Range("C9").Value = 1782
